$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-10) represent weekly price reports for the same
# market/product. This edit re-sorts them by date (column D, "Fecha")
# in descending order, carrying along the related volume/price columns
# M (Volumen), N (Precio minimo), O (Precio maximo), P (Precio promedio
# ponderado) and S (Precio $/Kg). Other columns are identical across all
# rows, so only D, M, N, O, P, S need to move.

# Mapping: target row -> source row (by original row order) after sorting
# descending on column D.
$rowMap = @{
    2  = 9
    3  = 4
    4  = 3
    5  = 2
    6  = 5
    7  = 7
    8  = 8
    9  = 10
    10 = 6
}

# Snapshot the original values for columns D, M, N, O, P, S before
# overwriting anything.
$cols = @("D", "M", "N", "O", "P", "S")
$orig = @{}
foreach ($r in 2..10) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value2 = $orig[$sourceRow][$c]
    }
}
